$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 859.75
$ws.Range("J17").Value = 859.75
$ws.Range("L17").Value = 2579.25
$ws.Range("N17").Value = -2915.25
$ws.Range("H137").Value = 1438.2
$ws.Range("I137").Value = 1264.6666
$ws.Range("K137").Value = 3793.9998
$ws.Range("M137").Value = -1243.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3463.9
$ws.Range("I32").Value = 2151.4634
$ws.Range("K32").Value = 2151.4634
$ws.Range("M32").Value = -1864.4634
$ws.Range("H122").Value = 14947.833
$ws.Range("I122").Value = 10525.5
$ws.Range("J122").Value = 30426
$ws.Range("K122").Value = 31576.5
$ws.Range("L122").Value = 91278
$ws.Range("M122").Value = -29126.5
$ws.Range("N122").Value = -96178

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3059.3696
$ws.Range("I31").Value = 2140.9412
$ws.Range("J31").Value = 5661.5835
$ws.Range("K31").Value = 2140.9412
$ws.Range("L31").Value = 5661.5835
$ws.Range("M31").Value = -1845.9412
$ws.Range("N31").Value = -6251.5835
$ws.Range("H34").Value = 3059.3696
$ws.Range("I34").Value = 2140.9412
$ws.Range("J34").Value = 5661.5835
$ws.Range("K34").Value = 2140.9412
$ws.Range("L34").Value = 5661.5835
$ws.Range("M34").Value = -1938.9412
$ws.Range("N34").Value = -6065.5835

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 62561.812
$ws.Range("I33").Value = 68.166664
$ws.Range("J33").Value = 250042.75
$ws.Range("K33").Value = 408.999984
$ws.Range("L33").Value = 1500256.5
$ws.Range("M33").Value = -125.999984
$ws.Range("N33").Value = -1500822.5
$ws.Range("H92").Value = 590
$ws.Range("J92").Value = 580
$ws.Range("L92").Value = 1740
$ws.Range("N92").Value = -4236

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 7396.8
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 7396.8
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 7396.8
$ws.Range("N7").Value = -7620.8
$ws.Range("M7").ClearContents()
$ws.Range("H8").Value = 7396.8
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 7396.8
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 7396.8
$ws.Range("N8").Value = -7674.8
$ws.Range("M8").ClearContents()
$ws.Range("H24").Value = 24075
$ws.Range("I24").Value = 9006
$ws.Range("J24").Value = 25510.143
$ws.Range("K24").Value = 9006
$ws.Range("L24").Value = 25510.143
$ws.Range("M24").Value = -8833
$ws.Range("N24").Value = -25856.143
$ws.Range("H26").Value = 41706.332
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 41706.332
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 41706.332
$ws.Range("N26").Value = -42266.332
$ws.Range("M26").ClearContents()
$ws.Range("H29").Value = 19996
$ws.Range("I29").Value = 19996
$ws.Range("K29").Value = 19996
$ws.Range("M29").Value = -19706
$ws.Range("H50").Value = 41706.332
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 41706.332
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 41706.332
$ws.Range("N50").Value = -42702.332
$ws.Range("M50").ClearContents()
$ws.Range("H122").Value = 254991.75
$ws.Range("J122").Value = 338989.34
$ws.Range("L122").Value = 1016968.02
$ws.Range("N122").Value = -1021868.02

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3609.842
$ws.Range("I7").Value = 3370.4614
$ws.Range("J7").Value = 4128.5
$ws.Range("K7").Value = 3370.4614
$ws.Range("L7").Value = 4128.5
$ws.Range("M7").Value = -3258.4614
$ws.Range("N7").Value = -4352.5
$ws.Range("H16").Value = 1500
$ws.Range("I16").Value = 1500
$ws.Range("K16").Value = 1500
$ws.Range("M16").Value = -1330
$ws.Range("H22").Value = 37500
$ws.Range("I22").Value = 50000
$ws.Range("K22").Value = 50000
$ws.Range("M22").Value = -49705
$ws.Range("H27").Value = 37500
$ws.Range("I27").Value = 50000
$ws.Range("K27").Value = 50000
$ws.Range("M27").Value = -49893
$ws.Range("H40").Value = 4973.5
$ws.Range("I40").Value = 4777.6
$ws.Range("K40").Value = 4777.6
$ws.Range("M40").Value = -4641.6
$ws.Range("H46").Value = 2140
$ws.Range("I46").Value = 3000
$ws.Range("K46").Value = 3000
$ws.Range("M46").Value = -2812
$ws.Range("H82").Value = 85466.086
$ws.Range("I82").Value = 1843.7778
$ws.Range("J82").Value = 336333
$ws.Range("K82").Value = 1843.7778
$ws.Range("L82").Value = 336333
$ws.Range("M82").Value = -1482.7778
$ws.Range("N82").Value = -337055
$ws.Range("H85").Value = 85466.086
$ws.Range("I85").Value = 1843.7778
$ws.Range("J85").Value = 336333
$ws.Range("K85").Value = 1843.7778
$ws.Range("L85").Value = 336333
$ws.Range("M85").Value = -595.7778000000001
$ws.Range("N85").Value = -338829
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H124").Value = 63429
$ws.Range("J124").Value = 63429
$ws.Range("L124").Value = 63429
$ws.Range("N124").Value = -73249
$ws.Range("H126").Value = 3609.842
$ws.Range("I126").Value = 3370.4614
$ws.Range("J126").Value = 4128.5
$ws.Range("K126").Value = 10111.3842
$ws.Range("L126").Value = 12385.5
$ws.Range("M126").Value = -7641.3842
$ws.Range("N126").Value = -17325.5
$ws.Range("H132").Value = 501752
$ws.Range("I132").Value = 501752
$ws.Range("K132").Value = 1505256
$ws.Range("M132").Value = -1502726

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H96").Value = 1735.4286
$ws.Range("I96").Value = 1735.4286
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1735.4286
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = -362.4286
$ws.Range("M96").ClearContents()
$ws.Range("H107").Value = 537.2
$ws.Range("I107").Value = 624
$ws.Range("K107").Value = 1872
$ws.Range("M107").Value = 48
$ws.Range("H122").Value = 4022.1765
$ws.Range("I122").Value = 4089.3635
$ws.Range("K122").Value = 12268.0905
$ws.Range("M122").Value = -9818.0905
